# Add a new "2022-Q1" sheet (fund-level holdings detail) positioned right
# before the "总计" (Total) summary sheet, matching how each quarter's
# detail sheet already looks (2021-Q4, 2021-Q3, ...).

$wb = $excel.ActiveWorkbook
$totalSheetRef = $wb.Worksheets.Item("总计")

$newSheet = $wb.Worksheets.Add($totalSheetRef)
$newSheet.Name = "2022-Q1"

# Worksheet handles track by position, not identity, so re-resolve "总计"
# by name now that the new sheet has shifted everything after it by one.
$totalSheet = $wb.Worksheets.Item("总计")

function Set-HeaderCell($range, $text) {
    $range.Value = $text
    $range.Font.Bold = $true
    $range.HorizontalAlignment = -4108
    $range.VerticalAlignment = -4160
    $range.Borders.LineStyle = 1
}

function Set-TextValue($range, $text) {
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.Style = "Normal"
}

function Set-IndexCell($range, $num) {
    $range.Value = $num
    $range.Font.Bold = $true
    $range.HorizontalAlignment = -4108
    $range.VerticalAlignment = -4160
    $range.Borders.LineStyle = 1
}

# ---- header row (bold, boxed, centered/top like the other quarter sheets) ----
Set-HeaderCell $newSheet.Range("B1") "基金代码"
Set-HeaderCell $newSheet.Range("C1") "基金名称"
Set-HeaderCell $newSheet.Range("D1") "基金规模"
Set-HeaderCell $newSheet.Range("E1") "股票总仓位"
Set-HeaderCell $newSheet.Range("F1") "仓位占比"
Set-HeaderCell $newSheet.Range("G1") "持有市值(亿元)"
Set-HeaderCell $newSheet.Range("H1") "仓位排名"

# ---- data rows ----
Set-IndexCell $newSheet.Range("A2") 0
Set-TextValue $newSheet.Range("B2") "513090"
Set-TextValue $newSheet.Range("C2") "易方达中证香港证券投资主题ETF"
Set-TextValue $newSheet.Range("D2") "11.07"
Set-TextValue $newSheet.Range("E2") "96.47"
Set-TextValue $newSheet.Range("F2") "6.97"
Set-TextValue $newSheet.Range("G2") "0.7716"
$newSheet.Range("H2").Value = 5

Set-IndexCell $newSheet.Range("A3") 1
Set-TextValue $newSheet.Range("B3") "002860"
Set-TextValue $newSheet.Range("C3") "前海开源沪港深新机遇灵活配置混合"
Set-TextValue $newSheet.Range("D3") "0.01"
Set-TextValue $newSheet.Range("E3") "83.26"
Set-TextValue $newSheet.Range("F3") "6.36"
Set-TextValue $newSheet.Range("G3") "0.0006"
$newSheet.Range("H3").Value = 8

# ---- update the "总计" (Total) summary sheet: insert a new top row for
#      2022-Q1 and push the existing quarters down by one ----
$totalSheet.Rows("2").Insert()
# Inserting a row copies the formatting of the row above (the bold header),
# so reset the fresh row back to the plain/default look used by the other
# data rows before filling it in.
$totalSheet.Range("A2:D2").Style = "Normal"

Set-IndexCell $totalSheet.Range("A2") 0
$totalSheet.Range("B2").Value = "2022-Q1"
$totalSheet.Range("C2").Value = 2
$totalSheet.Range("D2").Value = 0.77

# Re-number the index column (A) for the rows that got pushed down so it
# stays a contiguous 0-based sequence, matching the existing convention.
for ($r = 3; $r -le 7; $r++) {
    $totalSheet.Range("A$r").Value = $r - 2
}

Write-Output "edit complete"
